$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 237
$ws.Range("F4").Value = 4849
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 1209
$ws.Range("F13").Value = 117
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 196
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 116
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 6401
$ws.Range("F22").Value = 0
$ws.Range("F25").Value = 545
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 28
$ws.Range("F31").Value = 2597
$ws.Range("F33").Value = 535
$ws.Range("F34").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 80

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 0

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 237
$ws.Range("F4").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 234
$ws.Range("F13").Value = 1209
$ws.Range("F14").Value = 117
$ws.Range("F15").Value = 196
$ws.Range("F16").Value = 85
$ws.Range("F18").Value = 156
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 6401
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 4001
$ws.Range("F28").Value = 0
$ws.Range("F30").Value = 28
$ws.Range("F31").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("F36").Value = 319
$ws.Range("F37").Value = 379
$ws.Range("F38").Value = 0
$ws.Range("F39").Value = 13
$ws.Range("F40").Value = 1573
$ws.Range("F41").Value = 979
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 60
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 0
